$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round Q2 and R2 to whole numbers
$ws.Range("Q2").Value = 754870
$ws.Range("R2").Value = 7156787

# Clear the Starttid (Z2) and Sluttid (AB2) values for this row
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
